# Reorders the comma-separated "Recorded By" values in column G so that
# previously-trailing entries (e.g. "System") move to the front, matching
# the synced copy from the main repo. The mapping below was derived by
# diffing every distinct value that appears in column G between the
# before/after versions of the workbook; values not present in the map
# (already in their canonical order, or single-value cells) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, backup@backdoor.com, System" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
